$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, preserving exact
# formatting (leading/trailing zeros, multi-dot thousand separators, etc.)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.797.53"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.493.28"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.18"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.43"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.515.19"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.938.92"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.11"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.754.79"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.502.74"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.24"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.90"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.85"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.56"
$ws.Range("E24").Value = "  +5.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.418"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.74"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0771"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.92"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.41"
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.08"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.76"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.28"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.72"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.88"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.39"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.30"
$ws.Range("E51").Value = "  -1.41%  "
